$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting row (07/11/2023) - copy the date formatting from the row above
# and then fill in the actual values for the new row.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = Get-Date -Year 2023 -Month 11 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("B10").Value = "Yes"
$ws.Range("C10").Value = "Yes"
$ws.Range("D10").Value = "Yes"
$ws.Range("E10").Value = "Yes"

$ws.Range("F10").Select()
